{"js": "// Applies the cover-letter revision described by the diff:\n//  - adds a dateline above the From/To header\n//  - tightens the From/To header block\n//  - rewrites/expands several body paragraphs\n//  - splits the old 3rd body paragraph into three focused paragraphs\n//\n// The whole letter lives in ONE <w:p>/<w:r>, with <w:br/> soft line\n// breaks marking \"paragraphs\" (there are no real paragraph marks), so\n// every edit below is a body-level text search & replace, matching the\n// canonical OOXML shape (alternating <w:t>/<w:br/>) exactly.\n\nconst body = context.document.body;\n\n// 1) Insert the dateline above everything else.\nbody.insertText(\"May 22nd, 2023\\u000b\\u000b\", Word.InsertLocation.start);\nawait context.sync();\n\n// 2) \"From / To\" header: the break after \"From: ...\" becomes a double\n//    break, and \"Hiring Manager, \" is dropped from the \"To:\" line.\nlet r1 = body.search(\"From: Austing Dong\\u000bTo: Hiring Manager, University of Waterloo - Faculty of Science\", { matchCase: true });\nr1.load(\"items\");\nawait context.sync();\nif (r1.items.length !== 1) {\n  throw new Error(\"edit 1: expected exactly 1 match, got \" + r1.items.length);\n}\nr1.items[0].insertText(\"From: Austing Dong\\u000b\\u000bTo: University of Waterloo - Faculty of Science\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Opening paragraph: add the closing sentence about highlights.\nlet r2 = body.search(\"I am writing to express my strong interest in applying for the position of Software Developer for Mac OS (Co-op) at the University of Waterloo - Faculty of Science. As a Computer Science undergraduate student at the University of Waterloo, I strongly believe that my technical competencies and academic background are closely in line with the job requirements.\", { matchCase: true });\nr2.load(\"items\");\nawait context.sync();\nif (r2.items.length !== 1) {\n  throw new Error(\"edit 2: expected exactly 1 match, got \" + r2.items.length);\n}\nr2.items[0].insertText(\"I am writing to express my strong interest in applying for the position of Software Developer for Mac OS (Co-op) at the University of Waterloo - Faculty of Science. As a Computer Science undergraduate student at the University of Waterloo, I strongly believe that my technical competencies and academic background are closely in line with the job requirements. I would like to highlight the following for your consideration:\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Origin-story paragraph: rewritten/expanded.\nlet r3 = body.search(\"I have always been passionate about computer science and developing applications since middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest. I have excelled in relevant courses in the beginning of my university studies and found solving business challenges through programming fascinating because it gives me a sense of accomplishment. Such deep interest in programming and technology has motivated me to deep dive in related fields such as software development, quality assurance, and machine learning.\", { matchCase: true });\nr3.load(\"items\");\nawait context.sync();\nif (r3.items.length !== 1) {\n  throw new Error(\"edit 3: expected exactly 1 match, got \" + r3.items.length);\n}\nr3.items[0].insertText(\"My passion for computer science and application development began in middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest. This interest continued to grow as I excelled in relevant courses during the early stages of my university studies. I find solving business challenges through programming fascinating, as it provides me with a sense of accomplishment. This deep interest in programming and technology has motivated me to explore related fields such as software development, quality assurance, and machine learning.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 5) Co-op paragraph: trim the GitHub sentence and move the project\n//    write-up out into two new paragraphs of their own.\nlet r4 = body.search(\"The computer science co-op program at the University of Waterloo offered me a unique opportunity to take on both programming and logical courses. Through working on a massive number of technical projects and assignments from hackathons and school activities, I gained hands-on experience in fields including but not limited to object-oriented programming, web application development, artificial intelligence, algorithm design, and data abstraction. Such projects can be viewed on my GitHub: https://github.com/AustingDong. One of the biggest projects I led and built was implementing an application that uses AI to extract keywords from articles containing scientific or technical information which helps users quickly locate their desired items based on keywords' weight. This application can be used to quickly get all the important items and keywords from NASA Technical Report Server which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project details can be found here. Through understanding the project requirements, researching on coding algorithms, implementing the application with clean code, and incorporating white and black testing in a limited time, I cultivated working experience in developing a technical application in accordance with business requirements. Moreover, I have excellent documentation and communication skills through interpreting and explaining technical concepts to my teammates while working in a team environment.\", { matchCase: true });\nr4.load(\"items\");\nawait context.sync();\nif (r4.items.length !== 1) {\n  throw new Error(\"edit 4: expected exactly 1 match, got \" + r4.items.length);\n}\nr4.items[0].insertText(\"The computer science co-op program at the University of Waterloo has offered me a unique opportunity to take on both programming and logical courses. Through working on numerous technical projects and assignments from hackathons and school activities, I have gained hands-on experience in fields including, but not limited to, object-oriented programming, web application development, artificial intelligence, algorithm design, and data abstraction. My projects can be viewed on my GitHub: https://github.com/AustingDong.\\u000b\\u000bOne of the most significant projects I led and built was implementing an application that uses AI to extract keywords from articles containing scientific or technical information. This application helps users quickly locate their desired items based on keyword weight and can be used to efficiently retrieve important items and keywords from the NASA Technical Report Server, which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project details can be found here.\\u000b\\u000bThrough understanding project requirements, researching coding algorithms, implementing the application with clean code, and incorporating white and black testing within a limited time frame, I have cultivated valuable experience in developing technical applications in accordance with business requirements. Moreover, I possess excellent documentation and communication skills, which have been honed through interpreting and explaining technical concepts to my teammates while working in a team environment. I am confident that my academic and project background has prepared me well to make valuable contributions to a workplace environment.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 6) Closing-ask paragraph: drop the \"I am confident...\" sentence\n//    (it reappears, reworded, at the end of the project write-up).\nlet r5 = body.search(\"I am extremely interested in advancing my career and contributing my skills to the University of Waterloo - Faculty of Science. I am confident that my academic and project background has prepared me well, and now is the time I can make valuable contributions to a workplace environment. I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any questions or require additional information.\", { matchCase: true });\nr5.load(\"items\");\nawait context.sync();\nif (r5.items.length !== 1) {\n  throw new Error(\"edit 5: expected exactly 1 match, got \" + r5.items.length);\n}\nr5.items[0].insertText(\"I am extremely interested in advancing my career and contributing my skills to the University of Waterloo - Faculty of Science. I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any questions or require additional information.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Applies the cover-letter revision described by the diff:\n#  - adds a dateline above the From/To header\n#  - tightens the From/To header block\n#  - rewrites/expands several body paragraphs\n#  - splits the old 3rd body paragraph into three focused paragraphs\n#\n# The whole letter lives in ONE <w:p>/<w:r>, with <w:br/> soft line\n# breaks marking \"paragraphs\" (there are no real paragraph marks), so\n# every edit below is a Find/Replace over $d.Content, matching the\n# canonical OOXML shape (alternating <w:t>/<w:br/>) exactly.\n# [char]11 is a Word \"manual line break\" (<w:br/>), i.e. chr(11)/VT.\n\n$d = $word.ActiveDocument\n\n# 1) Insert the dateline above everything else.\n$d.Range(0, 0).InsertBefore(\"May 22nd, 2023\" + [char]11 + [char]11)\n\n# 2) \"From / To\" header: the break after \"From: ...\" becomes a double\n#    break, and \"Hiring Manager, \" is dropped from the \"To:\" line.\n$find1 = $d.Content.Find\n$find1.Execute(\"From: Austing Dong\" + [char]11 + \"To: Hiring Manager, University of Waterloo - Faculty of Science\", $true, $false, $false, $false, $false, $true, 1, $false, \"From: Austing Dong\" + [char]11 + [char]11 + \"To: University of Waterloo - Faculty of Science\", 2) | Out-Null\n\n# 3) Opening paragraph: add the closing sentence about highlights.\n$find2 = $d.Content.Find\n$find2.Execute(\"I am writing to express my strong interest in applying for the position of Software Developer for Mac OS (Co-op) at the University of Waterloo - Faculty of Science. As a Computer Science undergraduate student at the University of Waterloo, I strongly believe that my technical competencies and academic background are closely in line with the job requirements.\", $true, $false, $false, $false, $false, $true, 1, $false, \"I am writing to express my strong interest in applying for the position of Software Developer for Mac OS (Co-op) at the University of Waterloo - Faculty of Science. As a Computer Science undergraduate student at the University of Waterloo, I strongly believe that my technical competencies and academic background are closely in line with the job requirements. I would like to highlight the following for your consideration:\", 2) | Out-Null\n\n# 4) Origin-story paragraph: rewritten/expanded.\n$find3 = $d.Content.Find\n$find3.Execute(\"I have always been passionate about computer science and developing applications since middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest. I have excelled in relevant courses in the beginning of my university studies and found solving business challenges through programming fascinating because it gives me a sense of accomplishment. Such deep interest in programming and technology has motivated me to deep dive in related fields such as software development, quality assurance, and machine learning.\", $true, $false, $false, $false, $false, $true, 1, $false, \"My passion for computer science and application development began in middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest. This interest continued to grow as I excelled in relevant courses during the early stages of my university studies. I find solving business challenges through programming fascinating, as it provides me with a sense of accomplishment. This deep interest in programming and technology has motivated me to explore related fields such as software development, quality assurance, and machine learning.\", 2) | Out-Null\n\n# 5) Co-op paragraph: trim the GitHub sentence and move the project\n#    write-up out into two new paragraphs of their own.\n$find4 = $d.Content.Find\n$find4.Execute(\"The computer science co-op program at the University of Waterloo offered me a unique opportunity to take on both programming and logical courses. Through working on a massive number of technical projects and assignments from hackathons and school activities, I gained hands-on experience in fields including but not limited to object-oriented programming, web application development, artificial intelligence, algorithm design, and data abstraction. Such projects can be viewed on my GitHub: https://github.com/AustingDong. One of the biggest projects I led and built was implementing an application that uses AI to extract keywords from articles containing scientific or technical information which helps users quickly locate their desired items based on keywords' weight. This application can be used to quickly get all the important items and keywords from NASA Technical Report Server which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project details can be found here. Through understanding the project requirements, researching on coding algorithms, implementing the application with clean code, and incorporating white and black testing in a limited time, I cultivated working experience in developing a technical application in accordance with business requirements. Moreover, I have excellent documentation and communication skills through interpreting and explaining technical concepts to my teammates while working in a team environment.\", $true, $false, $false, $false, $false, $true, 1, $false, \"The computer science co-op program at the University of Waterloo has offered me a unique opportunity to take on both programming and logical courses. Through working on numerous technical projects and assignments from hackathons and school activities, I have gained hands-on experience in fields including, but not limited to, object-oriented programming, web application development, artificial intelligence, algorithm design, and data abstraction. My projects can be viewed on my GitHub: https://github.com/AustingDong.\" + [char]11 + [char]11 + \"One of the most significant projects I led and built was implementing an application that uses AI to extract keywords from articles containing scientific or technical information. This application helps users quickly locate their desired items based on keyword weight and can be used to efficiently retrieve important items and keywords from the NASA Technical Report Server, which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project details can be found here.\" + [char]11 + [char]11 + \"Through understanding project requirements, researching coding algorithms, implementing the application with clean code, and incorporating white and black testing within a limited time frame, I have cultivated valuable experience in developing technical applications in accordance with business requirements. Moreover, I possess excellent documentation and communication skills, which have been honed through interpreting and explaining technical concepts to my teammates while working in a team environment. I am confident that my academic and project background has prepared me well to make valuable contributions to a workplace environment.\", 2) | Out-Null\n\n# 6) Closing-ask paragraph: drop the \"I am confident...\" sentence\n#    (it reappears, reworded, at the end of the project write-up).\n$find5 = $d.Content.Find\n$find5.Execute(\"I am extremely interested in advancing my career and contributing my skills to the University of Waterloo - Faculty of Science. I am confident that my academic and project background has prepared me well, and now is the time I can make valuable contributions to a workplace environment. I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any questions or require additional information.\", $true, $false, $false, $false, $false, $true, 1, $false, \"I am extremely interested in advancing my career and contributing my skills to the University of Waterloo - Faculty of Science. I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any questions or require additional information.\", 2) | Out-Null\n"}
